# Apply the simulator full-month coverage edit:
#  - Fix the client names on the "Weekly Timesheet" sheet (and the mirrored
#    "Jason Schema" log sheet) which had drifted out of sync with the roster.
#  - Persist the simulator's computed Rate/Total figures (previously left at
#    0 because the log wasn't being flushed) on both sheets.
#  - Fix the employee id on the "Jason Schema" sheet.

$wb = $excel.ActiveWorkbook

$wsTime = $wb.Worksheets.Item("Weekly Timesheet")
$wsLog  = $wb.Worksheets.Item("Jason Schema")

# --- Client name corrections -------------------------------------------------
# "Weekly Timesheet": column B holds the client name, rows 2-6 = 2026-01-12..16
$wsTime.Range("B2").Value = "Hunter"
$wsTime.Range("B3").Value = "Tubergen"
$wsTime.Range("B4").Value = "Field"
$wsTime.Range("B5").Value = "Bottomley"
$wsTime.Range("B6").Value = "Zygmunt"

# "Jason Schema": column D mirrors the same client name, rows 2-6
$wsLog.Range("D2").Value = "Hunter"
$wsLog.Range("D3").Value = "Tubergen"
$wsLog.Range("D4").Value = "Field"
$wsLog.Range("D5").Value = "Bottomley"
$wsLog.Range("D6").Value = "Zygmunt"

# --- Employee ID correction ("Jason Schema" column B, rows 2-6) -------------
for ($r = 2; $r -le 6; $r++) {
    $wsLog.Cells.Item($r, 2).Value = "emp_yde33znx"
}

# --- Persist simulator-computed Rate / Total figures ------------------------
# "Weekly Timesheet": E = Rate, F = Total, for the daily rows 2-6
for ($r = 2; $r -le 6; $r++) {
    $wsTime.Cells.Item($r, 5).Value = 92
    $wsTime.Cells.Item($r, 6).Value = 736
}

# Roll the totals up through the summary rows (SUBTOTAL / HOURLY SUBTOTAL / GRAND TOTAL)
$wsTime.Range("F8").Value = 3680
$wsTime.Range("F11").Value = 3680
$wsTime.Range("F13").Value = 3680

# "Jason Schema": F = Rate, G = Total, for the daily rows 2-6
for ($r = 2; $r -le 6; $r++) {
    $wsLog.Cells.Item($r, 6).Value = 92
    $wsLog.Cells.Item($r, 7).Value = 736
}
